# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 941
$ws1.Range("F10").Value = 777
$ws1.Range("F17").Value = 497
$ws1.Range("F23").Value = 1357
$ws1.Range("F29").Value = 337
$ws1.Range("F30").Value = 2015
$ws1.Range("F31").Value = 94
$ws1.Range("F32").Value = 60

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 941
$ws4.Range("F17").Value = 777
$ws4.Range("F29").Value = 497
$ws4.Range("F35").Value = 1357
$ws4.Range("F43").Value = 337
$ws4.Range("F44").Value = 2016
$ws4.Range("F45").Value = 94
$ws4.Range("F46").Value = 60
